$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("PyCHAMobs")

# --- Header row additions: I1=GUAIACOL, J1=HCHO, K1=NO, L1=NO2 ---
$ws2.Cells.Item(1, 9).Value = "GUAIACOL"
$ws2.Cells.Item(1, 10).Value = "HCHO"
$ws2.Cells.Item(1, 11).Value = "NO"
$ws2.Cells.Item(1, 12).Value = "NO2"

# Apply scientific number format to K and L columns (rows 2-92)
$ws2.Range("K2:L92").NumberFormat = "0.00E+00"

$ws2.Cells.Item(2, 9).Value = 10900000000000
$ws2.Cells.Item(2, 10).Value = 1758465644
$ws2.Cells.Item(2, 11).Value = 888000000000
$ws2.Cells.Item(2, 12).Value = 1050000000000

$ws2.Cells.Item(3, 9).Value = 10574000000000
$ws2.Cells.Item(3, 10).Value = 1758465644
$ws2.Cells.Item(3, 11).Value = 906000000000
$ws2.Cells.Item(3, 12).Value = 988000000000

$ws2.Cells.Item(4, 9).Value = 10248000000000
$ws2.Cells.Item(4, 10).Value = 1758465644
$ws2.Cells.Item(4, 11).Value = 930000000000
$ws2.Cells.Item(4, 12).Value = 934000000000

$ws2.Cells.Item(5, 9).Value = 9922000000000
$ws2.Cells.Item(5, 10).Value = 1758465644
$ws2.Cells.Item(5, 11).Value = 892000000000
$ws2.Cells.Item(5, 12).Value = 984000000000

$ws2.Cells.Item(6, 9).Value = 9596000000000
$ws2.Cells.Item(6, 10).Value = 1758465644
$ws2.Cells.Item(6, 11).Value = 831000000000
$ws2.Cells.Item(6, 12).Value = 1050000000000

$ws2.Cells.Item(7, 9).Value = 9270000000000
$ws2.Cells.Item(7, 10).Value = 1758465644
$ws2.Cells.Item(7, 11).Value = 754000000000
$ws2.Cells.Item(7, 12).Value = 1120000000000

$ws2.Cells.Item(8, 9).Value = 9016000000000
$ws2.Cells.Item(8, 10).Value = 9056967018
$ws2.Cells.Item(8, 11).Value = 692000000000
$ws2.Cells.Item(8, 12).Value = 1180000000000

$ws2.Cells.Item(9, 9).Value = 8762000000000
$ws2.Cells.Item(9, 10).Value = 16167990656
$ws2.Cells.Item(9, 11).Value = 651000000000
$ws2.Cells.Item(9, 12).Value = 1250000000000

$ws2.Cells.Item(10, 9).Value = 8508000000000
$ws2.Cells.Item(10, 10).Value = 24617300084
$ws2.Cells.Item(10, 11).Value = 609000000000
$ws2.Cells.Item(10, 12).Value = 1280000000000

$ws2.Cells.Item(11, 9).Value = 8254000000000
$ws2.Cells.Item(11, 10).Value = 33123658774
$ws2.Cells.Item(11, 11).Value = 555000000000
$ws2.Cells.Item(11, 12).Value = 1310000000000

$ws2.Cells.Item(12, 9).Value = 8000000000000
$ws2.Cells.Item(12, 10).Value = 41499412781
$ws2.Cells.Item(12, 11).Value = 518000000000
$ws2.Cells.Item(12, 12).Value = 1330000000000

$ws2.Cells.Item(13, 9).Value = 7774000000000
$ws2.Cells.Item(13, 10).Value = 50468705728
$ws2.Cells.Item(13, 11).Value = 504000000000
$ws2.Cells.Item(13, 12).Value = 1340000000000

$ws2.Cells.Item(14, 9).Value = 7548000000000
$ws2.Cells.Item(14, 10).Value = 60686968653
$ws2.Cells.Item(14, 11).Value = 473000000000
$ws2.Cells.Item(14, 12).Value = 1350000000000

$ws2.Cells.Item(15, 9).Value = 7322000000000
$ws2.Cells.Item(15, 10).Value = 70490484996
$ws2.Cells.Item(15, 11).Value = 437000000000
$ws2.Cells.Item(15, 12).Value = 1340000000000

$ws2.Cells.Item(16, 9).Value = 7096000000000
$ws2.Cells.Item(16, 10).Value = 79871811884
$ws2.Cells.Item(16, 11).Value = 413000000000
$ws2.Cells.Item(16, 12).Value = 1350000000000

$ws2.Cells.Item(17, 9).Value = 6870000000000
$ws2.Cells.Item(17, 10).Value = 91996741205
$ws2.Cells.Item(17, 11).Value = 391000000000
$ws2.Cells.Item(17, 12).Value = 1340000000000

$ws2.Cells.Item(18, 9).Value = 6670000000000
$ws2.Cells.Item(18, 10).Value = 102537000000
$ws2.Cells.Item(18, 11).Value = 369000000000
$ws2.Cells.Item(18, 12).Value = 1340000000000

$ws2.Cells.Item(19, 9).Value = 6470000000000
$ws2.Cells.Item(19, 10).Value = 112486000000
$ws2.Cells.Item(19, 11).Value = 353000000000
$ws2.Cells.Item(19, 12).Value = 1320000000000

$ws2.Cells.Item(20, 9).Value = 6270000000000
$ws2.Cells.Item(20, 10).Value = 126114000000
$ws2.Cells.Item(20, 11).Value = 349000000000
$ws2.Cells.Item(20, 12).Value = 1310000000000

$ws2.Cells.Item(21, 9).Value = 6070000000000
$ws2.Cells.Item(21, 10).Value = 146740000000
$ws2.Cells.Item(21, 11).Value = 348000000000
$ws2.Cells.Item(21, 12).Value = 1300000000000

$ws2.Cells.Item(22, 9).Value = 5870000000000
$ws2.Cells.Item(22, 10).Value = 158982000000
$ws2.Cells.Item(22, 11).Value = 344000000000
$ws2.Cells.Item(22, 12).Value = 1290000000000

$ws2.Cells.Item(23, 9).Value = 5696000000000
$ws2.Cells.Item(23, 10).Value = 165527000000
$ws2.Cells.Item(23, 11).Value = 339000000000
$ws2.Cells.Item(23, 12).Value = 1270000000000

$ws2.Cells.Item(24, 9).Value = 5522000000000
$ws2.Cells.Item(24, 10).Value = 188466000000
$ws2.Cells.Item(24, 11).Value = 329000000000
$ws2.Cells.Item(24, 12).Value = 1250000000000

$ws2.Cells.Item(25, 9).Value = 5348000000000
$ws2.Cells.Item(25, 10).Value = 201501000000
$ws2.Cells.Item(25, 11).Value = 318000000000
$ws2.Cells.Item(25, 12).Value = 1230000000000

$ws2.Cells.Item(26, 9).Value = 5174000000000
$ws2.Cells.Item(26, 10).Value = 206117000000
$ws2.Cells.Item(26, 11).Value = 308000000000
$ws2.Cells.Item(26, 12).Value = 1220000000000

$ws2.Cells.Item(27, 9).Value = 5000000000000
$ws2.Cells.Item(27, 10).Value = 213442000000
$ws2.Cells.Item(27, 11).Value = 295000000000
$ws2.Cells.Item(27, 12).Value = 1220000000000

$ws2.Cells.Item(28, 9).Value = 4872000000000
$ws2.Cells.Item(28, 10).Value = 220717000000
$ws2.Cells.Item(28, 11).Value = 281000000000
$ws2.Cells.Item(28, 12).Value = 1210000000000

$ws2.Cells.Item(29, 9).Value = 4744000000000
$ws2.Cells.Item(29, 10).Value = 228777000000
$ws2.Cells.Item(29, 11).Value = 270000000000
$ws2.Cells.Item(29, 12).Value = 1170000000000

$ws2.Cells.Item(30, 9).Value = 4616000000000
$ws2.Cells.Item(30, 10).Value = 237628000000
$ws2.Cells.Item(30, 11).Value = 260000000000
$ws2.Cells.Item(30, 12).Value = 1150000000000

$ws2.Cells.Item(31, 9).Value = 4488000000000
$ws2.Cells.Item(31, 10).Value = 255626000000
$ws2.Cells.Item(31, 11).Value = 246000000000
$ws2.Cells.Item(31, 12).Value = 1130000000000

$ws2.Cells.Item(32, 9).Value = 4360000000000
$ws2.Cells.Item(32, 10).Value = 258358000000
$ws2.Cells.Item(32, 11).Value = 238000000000
$ws2.Cells.Item(32, 12).Value = 1110000000000

$ws2.Cells.Item(33, 9).Value = 4256000000000
$ws2.Cells.Item(33, 10).Value = 259209000000
$ws2.Cells.Item(33, 11).Value = 230000000000
$ws2.Cells.Item(33, 12).Value = 1090000000000

$ws2.Cells.Item(34, 9).Value = 4152000000000
$ws2.Cells.Item(34, 10).Value = 261623000000
$ws2.Cells.Item(34, 11).Value = 220000000000
$ws2.Cells.Item(34, 12).Value = 1070000000000

$ws2.Cells.Item(35, 9).Value = 4048000000000
$ws2.Cells.Item(35, 10).Value = 268158000000
$ws2.Cells.Item(35, 11).Value = 215000000000
$ws2.Cells.Item(35, 12).Value = 1050000000000

$ws2.Cells.Item(36, 9).Value = 3944000000000
$ws2.Cells.Item(36, 10).Value = 282585000000
$ws2.Cells.Item(36, 11).Value = 208000000000
$ws2.Cells.Item(36, 12).Value = 1030000000000

$ws2.Cells.Item(37, 9).Value = 3840000000000
$ws2.Cells.Item(37, 10).Value = 298319000000
$ws2.Cells.Item(37, 11).Value = 201000000000
$ws2.Cells.Item(37, 12).Value = 1010000000000

$ws2.Cells.Item(38, 9).Value = 3750000000000
$ws2.Cells.Item(38, 10).Value = 301686000000
$ws2.Cells.Item(38, 11).Value = 194000000000
$ws2.Cells.Item(38, 12).Value = 991000000000

$ws2.Cells.Item(39, 9).Value = 3660000000000
$ws2.Cells.Item(39, 10).Value = 303747000000
$ws2.Cells.Item(39, 11).Value = 186000000000
$ws2.Cells.Item(39, 12).Value = 960000000000

$ws2.Cells.Item(40, 9).Value = 3480000000000
$ws2.Cells.Item(40, 10).Value = 326838000000
$ws2.Cells.Item(40, 11).Value = 173000000000
$ws2.Cells.Item(40, 12).Value = 934000000000

$ws2.Cells.Item(41, 9).Value = 3390000000000
$ws2.Cells.Item(41, 10).Value = 370110000000
$ws2.Cells.Item(41, 11).Value = 167000000000
$ws2.Cells.Item(41, 12).Value = 917000000000

$ws2.Cells.Item(42, 9).Value = 3314000000000
$ws2.Cells.Item(42, 10).Value = 383311000000
$ws2.Cells.Item(42, 11).Value = 159000000000
$ws2.Cells.Item(42, 12).Value = 899000000000

$ws2.Cells.Item(43, 9).Value = 3238000000000
$ws2.Cells.Item(43, 10).Value = 369750000000
$ws2.Cells.Item(43, 11).Value = 155000000000
$ws2.Cells.Item(43, 12).Value = 880000000000

$ws2.Cells.Item(44, 9).Value = 3162000000000
$ws2.Cells.Item(44, 10).Value = 364267000000
$ws2.Cells.Item(44, 11).Value = 148000000000
$ws2.Cells.Item(44, 12).Value = 863000000000

$ws2.Cells.Item(45, 9).Value = 3086000000000
$ws2.Cells.Item(45, 10).Value = 363907000000
$ws2.Cells.Item(45, 11).Value = 143000000000
$ws2.Cells.Item(45, 12).Value = 846000000000

$ws2.Cells.Item(46, 9).Value = 3010000000000
$ws2.Cells.Item(46, 10).Value = 369246000000
$ws2.Cells.Item(46, 11).Value = 138000000000
$ws2.Cells.Item(46, 12).Value = 826000000000

$ws2.Cells.Item(47, 9).Value = 2956666666666.6602
$ws2.Cells.Item(47, 10).Value = 387357000000
$ws2.Cells.Item(47, 11).Value = 137000000000
$ws2.Cells.Item(47, 12).Value = 809000000000

$ws2.Cells.Item(48, 9).Value = 2903333333333.3301
$ws2.Cells.Item(48, 10).Value = 401610000000
$ws2.Cells.Item(48, 11).Value = 132000000000
$ws2.Cells.Item(48, 12).Value = 792000000000

$ws2.Cells.Item(49, 9).Value = 2850000000000
$ws2.Cells.Item(49, 10).Value = 406518000000
$ws2.Cells.Item(49, 11).Value = 130000000000
$ws2.Cells.Item(49, 12).Value = 766000000000

$ws2.Cells.Item(50, 9).Value = 2796666666666.6602
$ws2.Cells.Item(50, 10).Value = 404686000000
$ws2.Cells.Item(50, 11).Value = 124000000000
$ws2.Cells.Item(50, 12).Value = 754000000000

$ws2.Cells.Item(51, 9).Value = 2743333333333.3301
$ws2.Cells.Item(51, 10).Value = 419110000000
$ws2.Cells.Item(51, 11).Value = 122000000000
$ws2.Cells.Item(51, 12).Value = 742000000000

$ws2.Cells.Item(52, 9).Value = 2690000000000
$ws2.Cells.Item(52, 10).Value = 429770000000
$ws2.Cells.Item(52, 11).Value = 119000000000
$ws2.Cells.Item(52, 12).Value = 732000000000

$ws2.Cells.Item(53, 9).Value = 2640000000000
$ws2.Cells.Item(53, 10).Value = 468560000000
$ws2.Cells.Item(53, 11).Value = 117000000000
$ws2.Cells.Item(53, 12).Value = 715000000000

$ws2.Cells.Item(54, 9).Value = 2590000000000
$ws2.Cells.Item(54, 10).Value = 469350000000
$ws2.Cells.Item(54, 11).Value = 113000000000
$ws2.Cells.Item(54, 12).Value = 701000000000

$ws2.Cells.Item(55, 9).Value = 2540000000000
$ws2.Cells.Item(55, 10).Value = 456237000000
$ws2.Cells.Item(55, 11).Value = 111000000000
$ws2.Cells.Item(55, 12).Value = 684000000000

$ws2.Cells.Item(56, 9).Value = 2490000000000
$ws2.Cells.Item(56, 10).Value = 437304000000
$ws2.Cells.Item(56, 11).Value = 105000000000
$ws2.Cells.Item(56, 12).Value = 673000000000

$ws2.Cells.Item(57, 9).Value = 2440000000000
$ws2.Cells.Item(57, 10).Value = 429154000000
$ws2.Cells.Item(57, 11).Value = 102000000000
$ws2.Cells.Item(57, 12).Value = 667000000000

$ws2.Cells.Item(58, 9).Value = 2396000000000
$ws2.Cells.Item(58, 10).Value = 460321000000
$ws2.Cells.Item(58, 11).Value = 99900000000
$ws2.Cells.Item(58, 12).Value = 667000000000

$ws2.Cells.Item(59, 9).Value = 2352000000000
$ws2.Cells.Item(59, 10).Value = 477961000000
$ws2.Cells.Item(59, 11).Value = 98200000000
$ws2.Cells.Item(59, 12).Value = 637000000000

$ws2.Cells.Item(60, 9).Value = 2308000000000
$ws2.Cells.Item(60, 10).Value = 483417000000
$ws2.Cells.Item(60, 11).Value = 96200000000
$ws2.Cells.Item(60, 12).Value = 618000000000

$ws2.Cells.Item(61, 9).Value = 2264000000000
$ws2.Cells.Item(61, 10).Value = 474513000000
$ws2.Cells.Item(61, 11).Value = 93500000000
$ws2.Cells.Item(61, 12).Value = 603000000000

$ws2.Cells.Item(62, 9).Value = 2220000000000
$ws2.Cells.Item(62, 10).Value = 456301000000
$ws2.Cells.Item(62, 11).Value = 91500000000
$ws2.Cells.Item(62, 12).Value = 580000000000

$ws2.Cells.Item(63, 9).Value = 2176000000000
$ws2.Cells.Item(63, 10).Value = 453540000000
$ws2.Cells.Item(63, 11).Value = 87300000000
$ws2.Cells.Item(63, 12).Value = 569000000000

$ws2.Cells.Item(64, 9).Value = 2132000000000
$ws2.Cells.Item(64, 10).Value = 674240000000
$ws2.Cells.Item(64, 11).Value = 85100000000
$ws2.Cells.Item(64, 12).Value = 561000000000

$ws2.Cells.Item(65, 9).Value = 2088000000000
$ws2.Cells.Item(65, 10).Value = 808212000000
$ws2.Cells.Item(65, 11).Value = 85100000000
$ws2.Cells.Item(65, 12).Value = 550000000000

$ws2.Cells.Item(66, 9).Value = 2044000000000
$ws2.Cells.Item(66, 10).Value = 653666000000
$ws2.Cells.Item(66, 11).Value = 83600000000
$ws2.Cells.Item(66, 12).Value = 540000000000

$ws2.Cells.Item(67, 9).Value = 2000000000000
$ws2.Cells.Item(67, 10).Value = 547203000000
$ws2.Cells.Item(67, 11).Value = 80200000000
$ws2.Cells.Item(67, 12).Value = 528000000000

$ws2.Cells.Item(68, 9).Value = 1972000000000
$ws2.Cells.Item(68, 10).Value = 572195000000
$ws2.Cells.Item(68, 11).Value = 77000000000
$ws2.Cells.Item(68, 12).Value = 518000000000

$ws2.Cells.Item(69, 9).Value = 1944000000000
$ws2.Cells.Item(69, 10).Value = 571346000000
$ws2.Cells.Item(69, 11).Value = 74500000000
$ws2.Cells.Item(69, 12).Value = 517000000000

$ws2.Cells.Item(70, 9).Value = 1916000000000
$ws2.Cells.Item(70, 10).Value = 539311000000
$ws2.Cells.Item(70, 11).Value = 74000000000
$ws2.Cells.Item(70, 12).Value = 499000000000

$ws2.Cells.Item(71, 9).Value = 1888000000000
$ws2.Cells.Item(71, 10).Value = 518566000000
$ws2.Cells.Item(71, 11).Value = 73300000000
$ws2.Cells.Item(71, 12).Value = 495000000000

$ws2.Cells.Item(72, 9).Value = 1860000000000
$ws2.Cells.Item(72, 10).Value = 539746000000
$ws2.Cells.Item(72, 11).Value = 72800000000
$ws2.Cells.Item(72, 12).Value = 491000000000

$ws2.Cells.Item(73, 9).Value = 1838000000000
$ws2.Cells.Item(73, 10).Value = 610235000000
$ws2.Cells.Item(73, 11).Value = 69600000000
$ws2.Cells.Item(73, 12).Value = 470000000000

$ws2.Cells.Item(74, 9).Value = 1816000000000
$ws2.Cells.Item(74, 10).Value = 587095000000
$ws2.Cells.Item(74, 11).Value = 69900000000
$ws2.Cells.Item(74, 12).Value = 459000000000

$ws2.Cells.Item(75, 9).Value = 1794000000000
$ws2.Cells.Item(75, 10).Value = 546138000000
$ws2.Cells.Item(75, 11).Value = 69400000000
$ws2.Cells.Item(75, 12).Value = 454000000000

$ws2.Cells.Item(76, 9).Value = 1772000000000
$ws2.Cells.Item(76, 10).Value = 549535000000
$ws2.Cells.Item(76, 11).Value = 71600000000
$ws2.Cells.Item(76, 12).Value = 447000000000

$ws2.Cells.Item(77, 9).Value = 1750000000000
$ws2.Cells.Item(77, 10).Value = 638143000000
$ws2.Cells.Item(77, 11).Value = 68100000000
$ws2.Cells.Item(77, 12).Value = 441000000000

$ws2.Cells.Item(78, 9).Value = 1734000000000
$ws2.Cells.Item(78, 10).Value = 661034000000
$ws2.Cells.Item(78, 11).Value = 68100000000
$ws2.Cells.Item(78, 12).Value = 433000000000

$ws2.Cells.Item(79, 9).Value = 1718000000000
$ws2.Cells.Item(79, 10).Value = 590304000000
$ws2.Cells.Item(79, 11).Value = 67400000000
$ws2.Cells.Item(79, 12).Value = 424000000000

$ws2.Cells.Item(80, 9).Value = 1702000000000
$ws2.Cells.Item(80, 10).Value = 548224000000
$ws2.Cells.Item(80, 11).Value = 64900000000
$ws2.Cells.Item(80, 12).Value = 415000000000

$ws2.Cells.Item(81, 9).Value = 1686000000000
$ws2.Cells.Item(81, 10).Value = 543938000000
$ws2.Cells.Item(81, 11).Value = 64700000000
$ws2.Cells.Item(81, 12).Value = 404000000000

$ws2.Cells.Item(82, 9).Value = 1670000000000
$ws2.Cells.Item(82, 10).Value = 560524000000
$ws2.Cells.Item(82, 11).Value = 62200000000
$ws2.Cells.Item(82, 12).Value = 395000000000

$ws2.Cells.Item(83, 9).Value = 1656000000000
$ws2.Cells.Item(83, 10).Value = 559543000000
$ws2.Cells.Item(83, 11).Value = 63200000000
$ws2.Cells.Item(83, 12).Value = 391000000000

$ws2.Cells.Item(84, 9).Value = 1642000000000
$ws2.Cells.Item(84, 10).Value = 546862000000
$ws2.Cells.Item(84, 11).Value = 62000000000
$ws2.Cells.Item(84, 12).Value = 388000000000

$ws2.Cells.Item(85, 9).Value = 1628000000000
$ws2.Cells.Item(85, 10).Value = 534270000000
$ws2.Cells.Item(85, 11).Value = 62000000000
$ws2.Cells.Item(85, 12).Value = 384000000000

$ws2.Cells.Item(86, 9).Value = 1614000000000
$ws2.Cells.Item(86, 10).Value = 556371000000
$ws2.Cells.Item(86, 11).Value = 59500000000
$ws2.Cells.Item(86, 12).Value = 377000000000

$ws2.Cells.Item(87, 9).Value = 1600000000000
$ws2.Cells.Item(87, 10).Value = 606189000000
$ws2.Cells.Item(87, 11).Value = 57600000000
$ws2.Cells.Item(87, 12).Value = 373000000000

$ws2.Cells.Item(88, 9).Value = 1580000000000
$ws2.Cells.Item(88, 10).Value = 612390000000
$ws2.Cells.Item(88, 11).Value = 57800000000
$ws2.Cells.Item(88, 12).Value = 371000000000

$ws2.Cells.Item(89, 9).Value = 1560000000000
$ws2.Cells.Item(89, 10).Value = 635276000000
$ws2.Cells.Item(89, 11).Value = 56800000000
$ws2.Cells.Item(89, 12).Value = 369000000000

$ws2.Cells.Item(90, 9).Value = 1540000000000
$ws2.Cells.Item(90, 10).Value = 629627000000
$ws2.Cells.Item(90, 11).Value = 56100000000
$ws2.Cells.Item(90, 12).Value = 361000000000

$ws2.Cells.Item(91, 9).Value = 1520000000000
$ws2.Cells.Item(91, 10).Value = 590209000000
$ws2.Cells.Item(91, 11).Value = 55400000000
$ws2.Cells.Item(91, 12).Value = 355000000000

$ws2.Cells.Item(92, 9).Value = 1500000000000
$ws2.Cells.Item(92, 10).Value = 571039000000
$ws2.Cells.Item(92, 11).Value = 55100000000
$ws2.Cells.Item(92, 12).Value = 352000000000

# --- Add new Sheet1 with data row ---
$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Cells.Item(1, 1).Value = 23777.923999999999
$ws3.Cells.Item(1, 2).Value = 5959813.6380000003
$ws3.Cells.Item(1, 3).Value = 6599434.875
$ws3.Cells.Item(1, 4).Value = 6468601.3140000002
$ws3.Cells.Item(1, 5).Value = 6179272.4950000001
$ws3.Cells.Item(1, 6).Value = 5916922.4759999998
$ws3.Cells.Item(1, 7).Value = 5640368.3119999999
$ws3.Cells.Item(1, 8).Value = 5396418.0250000004
$ws3.Cells.Item(1, 9).Value = 5195326.7220000001
$ws3.Cells.Item(1, 10).Value = 5054363.415
$ws3.Cells.Item(1, 11).Value = 4865105.2089999998
$ws3.Cells.Item(1, 12).Value = 4770989.9419999998
$ws3.Cells.Item(1, 13).Value = 4607930.45
$ws3.Cells.Item(1, 14).Value = 4441952.8720000004
$ws3.Cells.Item(1, 15).Value = 4329758.7869999995
$ws3.Cells.Item(1, 16).Value = 4213575.9790000003
$ws3.Cells.Item(1, 17).Value = 4142039.4010000001
$ws3.Cells.Item(1, 18).Value = 4062902.9789999998
$ws3.Cells.Item(1, 19).Value = 3985505.8450000002
$ws3.Cells.Item(1, 20).Value = 3885863.554
$ws3.Cells.Item(1, 21).Value = 3757664.656

# --- Restore selection/view state ---
$ws2.Activate()
$ws2.Range("M1:O1048576").Select()
